# Updates cached market-price / profit figures (columns H:N) on several
# per-sheet "Leve profit" tables, row by row, as refreshed by the market
# data scheduled runner. Columns are:
#   H=currentAveragePrice  I=currentAveragePriceNQ  J=currentAveragePriceHQ
#   K=LevePriceNQ  L=LevePriceHQ  M=LeveProfitNQ  N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58 (A Matter of Vital Importance / Mega-Potion of Vitality)
$ws.Range("H58").Value = 7474.9
$ws.Range("I58").Value = 124.833336
$ws.Range("J58").Value = 18500
$ws.Range("K58").Value = 374.500008
$ws.Range("L58").Value = 55500
$ws.Range("M58").Value = -224.500008
$ws.Range("N58").Value = -55800

# Row 98 (The Dotted Line / Enchanted Durium Ink)
$ws.Range("H98").Value = 11021.333
$ws.Range("I98").Value = 11024
$ws.Range("K98").Value = 11024
$ws.Range("M98").Value = -9526

# Row 100 (Asking for a Friend / Beetle Glue)
$ws.Range("H100").Value = 2858.7
$ws.Range("I100").Value = 2034.5
$ws.Range("J100").Value = 3408.1667
$ws.Range("K100").Value = 2034.5
$ws.Range("L100").Value = 3408.1667
$ws.Range("M100").Value = -1493.5
$ws.Range("N100").Value = -4490.1667

# Row 122 (Wishful Inking / Enchanted High Durium Ink)
$ws.Range("H122").Value = 11021.333
$ws.Range("I122").Value = 11024
$ws.Range("K122").Value = 33072
$ws.Range("M122").Value = -30622

# Row 135 (For Tired Minds / Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 345389.88
$ws.Range("I135").Value = 370885.4
$ws.Range("K135").Value = 3337968.6
$ws.Range("M135").Value = -3335433.6

# Row 137 (Cutting Edge of Culinary Quality / Magnesia Whetstone)
$ws.Range("H137").Value = 3805
$ws.Range("I137").Value = 10750
$ws.Range("K137").Value = 32250
$ws.Range("M137").Value = -29700

# Row 138 (All-night Crafting / Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 5321.394
$ws.Range("I138").Value = 1573.1364
$ws.Range("J138").Value = 12817.909
$ws.Range("K138").Value = 4719.4092
$ws.Range("L138").Value = 38453.727
$ws.Range("M138").Value = 420.5907999999999
$ws.Range("N138").Value = -48733.727

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Ain't Got No Ingots / Bronze Ingot)
$ws.Range("H2").Value = 4255.1665
$ws.Range("I2").Value = 2868.818
$ws.Range("K2").Value = 2868.818
$ws.Range("M2").Value = -2755.818

# Row 21 (Fashion Weak / Iron Cuirass)
$ws.Range("H21").Value = 2081.75
$ws.Range("I21").Value = 155
$ws.Range("J21").Value = 4008.5
$ws.Range("K21").Value = 155
$ws.Range("L21").Value = 4008.5
$ws.Range("M21").Value = 219
$ws.Range("N21").Value = -4756.5

# Row 24 (A Firm Hand / Iron Gauntlets)
$ws.Range("H24").Value = 42107.668
$ws.Range("J24").Value = 42107.668
$ws.Range("L24").Value = 42107.668
$ws.Range("N24").Value = -42855.668

# Row 36 (Hot for Teacher / Heavy Iron Armor)
$ws.Range("H36").Value = 8000
$ws.Range("I36").Value = 8000
$ws.Range("K36").Value = 8000
$ws.Range("M36").Value = -7654

# Row 45 (Hollow Hallmarks / Mythril Ingot)
$ws.Range("H45").Value = 9793.223
$ws.Range("I45").Value = 1750.3334
$ws.Range("J45").Value = 13814.667
$ws.Range("K45").Value = 1750.3334
$ws.Range("L45").Value = 13814.667
$ws.Range("M45").Value = -1373.3334
$ws.Range("N45").Value = -14568.667

# Row 74 (As the Bolt Flies / Titanium Nugget)
$ws.Range("H74").Value = 111157.734
$ws.Range("J74").Value = 5297.8
$ws.Range("L74").Value = 5297.8
$ws.Range("N74").Value = -7045.8

# Row 77 (Heavy Metal Banned (L) / Titanium Nugget)
$ws.Range("H77").Value = 111157.734
$ws.Range("J77").Value = 5297.8
$ws.Range("L77").Value = 26489
$ws.Range("N77").Value = -35225

# Row 100 (En Garde and on Guard / Doman Iron Gauntlets of Fending)
$ws.Range("H100").Value = 42107.668
$ws.Range("J100").Value = 42107.668
$ws.Range("L100").Value = 42107.668
$ws.Range("N100").Value = -44271.668

# Row 110 (Scheduled Maintenance / Deepgold Ingot)
$ws.Range("H110").Value = 15152527
$ws.Range("I110").Value = 830
$ws.Range("K110").Value = 830
$ws.Range("M110").Value = 1215

# Row 116 (No Scope / Titanbronze Ingot)
$ws.Range("H116").Value = 4255.1665
$ws.Range("I116").Value = 2868.818
$ws.Range("K116").Value = 2868.818
$ws.Range("M116").Value = -574.8180000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Hells Bells / Bronze Ingot)
$ws.Range("H3").Value = 4255.1665
$ws.Range("I3").Value = 2868.818
$ws.Range("K3").Value = 2868.818
$ws.Range("M3").Value = -2754.818

# Row 86 (Through Thick and Thin / Adamantite Nugget)
$ws.Range("H86").Value = 22729826
$ws.Range("I86").Value = 8335224
$ws.Range("K86").Value = 8335224
$ws.Range("M86").Value = -8334101

# Row 89 (Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget)
$ws.Range("H89").Value = 22729826
$ws.Range("I89").Value = 8335224
$ws.Range("K89").Value = 41676120
$ws.Range("M89").Value = -41670504

# Row 107 (The Gold Experience / Deepgold Nugget)
$ws.Range("H107").Value = 48917024
$ws.Range("I107").Value = 66179230
$ws.Range("J107").Value = 7433.6665
$ws.Range("K107").Value = 66179230
$ws.Range("L107").Value = 7433.6665
$ws.Range("M107").Value = -66177310
$ws.Range("N107").Value = -11273.6665

# Row 134 (Ruthenium Supremium / Ruthenium Ingot)
$ws.Range("H134").Value = 5620.4683
$ws.Range("I134").Value = 2327.0833
$ws.Range("J134").Value = 9057.044
$ws.Range("K134").Value = 6981.249899999999
$ws.Range("L134").Value = 27171.132
$ws.Range("M134").Value = -4446.249899999999
$ws.Range("N134").Value = -32241.132

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Driving Up the Wall / Elm Lumber)
$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50

# Row 134 (Wood You Be Quiet / Ceiba Lumber)
$ws.Range("H134").Value = 9193.772000000001
$ws.Range("I134").Value = 1887.3334
$ws.Range("K134").Value = 5662.0002
$ws.Range("M134").Value = -3127.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 43 (Sole Survivor / Baked Sole)
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 131 (The Mountain Steeped / Tsai tou Vounou)
$ws.Range("H131").Value = 1731.5834
$ws.Range("J131").Value = 3691.25
$ws.Range("L131").Value = 11073.75
$ws.Range("N131").Value = -21153.75

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (Copious Crystal Cannons / Manasilver Nugget)
$ws.Range("H113").Value = 7239.375
$ws.Range("I113").Value = 3749.625
$ws.Range("J113").Value = 8402.625
$ws.Range("K113").Value = 3749.625
$ws.Range("L113").Value = 8402.625
$ws.Range("M113").Value = -1579.625
$ws.Range("N113").Value = -12742.625

# Row 132 (On Board for Lar / Lar Ingot)
$ws.Range("H132").Value = 6788
$ws.Range("I132").Value = 3289.4285
$ws.Range("K132").Value = 9868.2855
$ws.Range("M132").Value = -7338.2855

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Tan Before the Ban / Leather)
$ws.Range("H7").Value = 6799.9287
$ws.Range("I7").Value = 4199
$ws.Range("K7").Value = 4199
$ws.Range("M7").Value = -4087

# Row 14 (Quelling Bloody Rumors / Hard Leather Shoes)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 22 (Skin off Their Backs / Aldgoat Leather)
$ws.Range("H22").Value = 3597.7273
$ws.Range("J22").Value = 4008.3333
$ws.Range("L22").Value = 4008.3333
$ws.Range("N22").Value = -4598.3333

# Row 27 (Fire and Hide / Aldgoat Leather)
$ws.Range("H27").Value = 3597.7273
$ws.Range("J27").Value = 4008.3333
$ws.Range("L27").Value = 4008.3333
$ws.Range("N27").Value = -4222.3333

# Row 61 (Spelling Me Softly / Raptor Leather)
$ws.Range("H61").Value = 5221.2
$ws.Range("I61").Value = 1672.9
$ws.Range("K61").Value = 1672.9
$ws.Range("M61").Value = -1470.9

# Row 93 (Hide to Go Seek / Gagana Leather)
$ws.Range("H93").Value = 3862.8386
$ws.Range("I93").Value = 4008.5
$ws.Range("J93").Value = 3707.4666
$ws.Range("K93").Value = 4008.5
$ws.Range("L93").Value = 3707.4666
$ws.Range("M93").Value = -2760.5
$ws.Range("N93").Value = -6203.4666

# Row 113 (Peace in Rest / Atrociraptor Leather)
$ws.Range("H113").Value = 5221.2
$ws.Range("I113").Value = 1672.9
$ws.Range("K113").Value = 1672.9
$ws.Range("M113").Value = 497.0999999999999

# Row 126 (Battered Books / Saiga Leather)
$ws.Range("H126").Value = 6799.9287
$ws.Range("I126").Value = 4199
$ws.Range("K126").Value = 12597
$ws.Range("M126").Value = -10127

# Row 132 (Tenets of Tanning / Silver Lobo Leather)
$ws.Range("H132").Value = 13521079
$ws.Range("I132").Value = 45457440
$ws.Range("J132").Value = 9541.77
$ws.Range("K132").Value = 136372320
$ws.Range("L132").Value = 28625.31
$ws.Range("M132").Value = -136369790
$ws.Range("N132").Value = -33685.31

# Row 136 (Respect for Br'aax / Br'aax Leather)
$ws.Range("H136").Value = 12184.637
$ws.Range("J136").Value = 12602.952
$ws.Range("L136").Value = 37808.856
$ws.Range("N136").Value = -42908.856

$ws = $wb.Worksheets.Item("WVR")
# Row 41 (Half Is the New Double / Linen Halfgloves)
$ws.Range("H41").Value = 18080.5
$ws.Range("J41").Value = 18080.5
$ws.Range("L41").Value = 18080.5
$ws.Range("N41").Value = -18860.5
